# Apply the workbook edit described by the diff:
#  - Change the "Mevcut_Sayisi" (D6) value for "Teknisyen/Tekniker" from 52 to 49
#    (the dependent "Fark" formula in E6 recalculates from 13 to 10 automatically)
#  - Update the sheet's current selection to the full used range A1:E11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value (dependent formula E6 = D6-C6 recalculates automatically)
$ws.Range("D6").Value = 49

# Select the whole used range A1:E11
$ws.Range("A1:E11").Select()
